$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$block1 = New-Object 'object[,]' 24,7
$block1[0,0] = 10.54512200968881
$block1[0,1] = 5.067007107086827
$block1[0,2] = 12.85918902245216
$block1[0,3] = 23.33821417090123
$block1[0,4] = 27.07302548811116
$block1[0,5] = 13.99644850093874
$block1[0,6] = 22.33477718897089
$block1[1,0] = 10.51544816081652
$block1[1,1] = 4.991445194677296
$block1[1,2] = 12.89147873705386
$block1[1,3] = 23.44567964130037
$block1[1,4] = 27.26150722008892
$block1[1,5] = 14.0647691829311
$block1[1,6] = 22.44588244739819
$block1[2,0] = 10.49914651585468
$block1[2,1] = 4.944037968147537
$block1[2,2] = 12.91375020834259
$block1[2,3] = 23.51920722058831
$block1[2,4] = 27.38839974761767
$block1[2,5] = 14.10937777439912
$block1[2,6] = 22.51981039603819
$block1[3,0] = 10.49299088974472
$block1[3,1] = 4.924478316026326
$block1[3,2] = 12.92344035219447
$block1[3,3] = 23.55105857170193
$block1[3,4] = 27.44289824160862
$block1[3,5] = 14.12822507385318
$block1[3,6] = 22.55136876479267
$block1[4,0] = 10.49199833139119
$block1[4,1] = 4.921216311750445
$block1[4,2] = 12.92508647912508
$block1[4,3] = 23.55646125558947
$block1[4,4] = 27.4521155671286
$block1[4,5] = 14.13139506219845
$block1[4,6] = 22.55669540524715
$block1[5,0] = 10.49906151898606
$block1[5,1] = 4.943775136627981
$block1[5,2] = 12.91387840640715
$block1[5,3] = 23.51962914504923
$block1[5,4] = 27.38912346630782
$block1[5,5] = 14.10962924688909
$block1[5,6] = 22.5202302091156
$block1[6,0] = 10.53449551989938
$block1[6,1] = 5.041171427037754
$block1[6,2] = 12.86981462577604
$block1[6,3] = 23.37369756675771
$block1[6,4] = 27.1356857908474
$block1[6,5] = 14.01945363048238
$block1[6,6] = 22.37189935861505
$block1[7,0] = 10.61895920834163
$block1[7,1] = 5.223560570659338
$block1[7,2] = 12.80283359432352
$block1[7,3] = 23.14775091428026
$block1[7,4] = 26.72810040732638
$block1[7,5] = 13.8637113786235
$block1[7,6] = 22.126472277103
$block1[8,0] = 10.68979308397474
$block1[8,1] = 5.35158044825927
$block1[8,2] = 12.76549639300785
$block1[8,3] = 23.01895273738118
$block1[8,4] = 26.4842708869674
$block1[8,5] = 13.76213247963652
$block1[8,6] = 21.97407103099207
$block1[9,0] = 10.72383958555173
$block1[9,1] = 5.408375644265742
$block1[9,2] = 12.75109390419187
$block1[9,3] = 22.96853621751432
$block1[9,4] = 26.38566184551242
$block1[9,5] = 13.71870856844906
$block1[9,6] = 21.9108453853041
$block1[10,0] = 10.7369866594004
$block1[10,1] = 5.429663711300769
$block1[10,2] = 12.74601165722439
$block1[10,3] = 22.95062781852732
$block1[10,4] = 26.35011069630536
$block1[10,5] = 13.70266550024858
$block1[10,6] = 21.88778473662091
$block1[11,0] = 10.73414401619657
$block1[11,1] = 5.425088867567583
$block1[11,2] = 12.74708967615148
$block1[11,3] = 22.95443196825839
$block1[11,4] = 26.35768735197786
$block1[11,5] = 13.70610284110212
$block1[11,6] = 21.89271199
$block1[12,0] = 10.72491615093285
$block1[12,1] = 5.41013148135055
$block1[12,2] = 12.75066833450146
$block1[12,3] = 22.96703912574566
$block1[12,4] = 26.38270104550956
$block1[12,5] = 13.71738066537431
$block1[12,6] = 21.90893048178051
$block1[13,0] = 10.71929670825087
$block1[13,1] = 5.400940786741059
$block1[13,2] = 12.75290877543115
$block1[13,3] = 22.9749156905495
$block1[13,4] = 26.39825634177343
$block1[13,5] = 13.72434083710206
$block1[13,6] = 21.91897969625655
$block1[14,0] = 10.68760410139537
$block1[14,1] = 5.347838722268869
$block1[14,2] = 12.76648962314378
$block1[14,3] = 23.02241273929462
$block1[14,4] = 26.49096462806443
$block1[14,5] = 13.76502638601647
$block1[14,6] = 21.9783261153909
$block1[15,0] = 10.66862336236314
$block1[15,1] = 5.314884544996044
$block1[15,2] = 12.77548270122127
$block1[15,3] = 23.0536499618238
$block1[15,4] = 26.5510052440698
$block1[15,5] = 13.7906990252437
$block1[15,6] = 22.01629906077073
$block1[16,0] = 10.65787830414274
$block1[16,1] = 5.295795139967472
$block1[16,2] = 12.7808983354799
$block1[16,3] = 23.07238576074387
$block1[16,4] = 26.5866965485852
$block1[16,5] = 13.80572735833427
$block1[16,6] = 22.03871424703148
$block1[17,0] = 10.65427001426601
$block1[17,1] = 5.289308970157665
$block1[17,2] = 12.78277370967932
$block1[17,3] = 23.07886120909161
$block1[17,4] = 26.59897921171361
$block1[17,5] = 13.81086070680624
$block1[17,6] = 22.04640213281293
$block1[18,0] = 10.67062613220951
$block1[18,1] = 5.318406635727599
$block1[18,2] = 12.77450021497543
$block1[18,3] = 23.05024506256555
$block1[18,4] = 26.54449387398175
$block1[18,5] = 13.78793900045159
$block1[18,6] = 22.0121973241923
$block1[19,0] = 10.72761975867386
$block1[19,1] = 5.414530863518808
$block1[19,2] = 12.74960710689249
$block1[19,3] = 22.96330392521635
$block1[19,4] = 26.37530518241086
$block1[19,5] = 13.71405722066009
$block1[19,6] = 21.90414276126566
$block1[20,0] = 10.76634734512672
$block1[20,1] = 5.4760718405204
$block1[20,2] = 12.73550443969792
$block1[20,3] = 22.91338302087174
$block1[20,4] = 26.2751727798355
$block1[20,5] = 13.66810645577795
$block1[20,6] = 21.83866294002819
$block1[21,0] = 10.74554506782201
$block1[21,1] = 5.443347266544266
$block1[21,2] = 12.74283299009679
$block1[21,3] = 22.93939292226957
$block1[21,4] = 26.32765341216918
$block1[21,5] = 13.6924175037928
$block1[21,6] = 21.8731391167079
$block1[22,0] = 10.66972015850409
$block1[22,1] = 5.316814744922548
$block1[22,2] = 12.77494363265602
$block1[22,3] = 23.05178199748743
$block1[22,4] = 26.54743401159164
$block1[22,5] = 13.78918597039445
$block1[22,6] = 22.01404990036634
$block1[23,0] = 10.59454073471941
$block1[23,1] = 5.175218596533263
$block1[23,2] = 12.81887043567937
$block1[23,3] = 23.20237615235872
$block1[23,4] = 26.82867555271945
$block1[23,5] = 13.90358724814936
$block1[23,6] = 22.18797976566563
$ws.Range("C2:I25").Value = $block1

$block2 = New-Object 'object[,]' 24,3
$block2[0,0] = 13.9089731205515
$block2[0,1] = 9.600923397982324
$block2[0,2] = 16.18939550126271
$block2[1,0] = 13.28499069818721
$block2[1,1] = 9.634373919280852
$block2[1,2] = 15.9240095914344
$block2[2,0] = 12.8857863024522
$block2[2,1] = 9.656202517535739
$block2[2,2] = 15.75980406871243
$block2[3,0] = 12.7192438723817
$block2[3,1] = 9.665422621307247
$block2[3,2] = 15.69264425397572
$block2[4,0] = 12.69136205945253
$block2[4,1] = 9.666973240791773
$block2[4,2] = 15.68147975250795
$block2[5,0] = 12.88355563913669
$block2[5,1] = 9.656325547510637
$block2[5,2] = 15.75889922342754
$block2[6,0] = 13.69726007810372
$block2[6,1] = 9.612189690187254
$block2[6,2] = 16.09819385468024
$block2[7,0] = 15.15877894590271
$block2[7,1] = 9.535854117136502
$block2[7,2] = 16.7505156503462
$block2[8,0] = 16.14351063311785
$block2[8,1] = 9.48597018382878
$block2[8,2] = 17.2178825334606
$block2[9,0] = 16.57100274041186
$block2[9,1] = 9.464617252095964
$block2[9,2] = 17.42712947751824
$block2[10,0] = 16.72986420052029
$block2[10,1] = 9.456723677108911
$block2[10,2] = 17.50582057721606
$block2[11,0] = 16.69578597307081
$block2[11,1] = 9.458415151114338
$block2[11,2] = 17.48889830101897
$block2[12,0] = 16.58413330022928
$block2[12,1] = 9.463963990323206
$block2[12,2] = 17.4336147371949
$block2[13,0] = 16.51534738397843
$block2[13,1] = 9.467387851024055
$block2[13,2] = 17.39967901529825
$block2[14,0] = 16.11515424350383
$block2[14,1] = 9.487392576438882
$block2[14,2] = 17.20413467995708
$block2[15,0] = 15.8643471602644
$block2[15,1] = 9.500007690830044
$block2[15,2] = 17.08326808546812
$block2[16,0] = 15.71816677655689
$block2[16,1] = 9.507389670467155
$block2[16,2] = 17.01343607438026
$block2[17,0] = 15.66834495058116
$block2[17,1] = 9.509910749088878
$block2[17,2] = 16.98974040937542
$block2[18,0] = 15.89124560609627
$block2[18,1] = 9.498651741777936
$block2[18,2] = 17.09616737414294
$block2[19,0] = 16.61701092381603
$block2[19,1] = 9.46232894469418
$block2[19,2] = 17.44986816909086
$block2[20,0] = 17.07370766954416
$block2[20,1] = 9.439710732690855
$block2[20,2] = 17.67782158652527
$block2[21,0] = 16.83159564745574
$block2[21,1] = 9.451680039545
$block2[21,2] = 17.5564722957784
$block2[22,0] = 15.87909100508336
$block2[22,1] = 9.49926436330237
$block2[22,2] = 17.09033667301048
$block2[23,0] = 14.77858934296353
$block2[23,1] = 9.555414221401143
$block2[23,2] = 16.57585743537433
$ws.Range("K2:M25").Value = $block2

$block3 = New-Object 'object[,]' 24,1
$block3[0,0] = 21.0559719092578
$block3[1,0] = 21.18104837733704
$block3[2,0] = 21.26334872735341
$block3[3,0] = 21.29826797289564
$block3[4,0] = 21.30414963359932
$block3[5,0] = 21.26381407149426
$block3[6,0] = 21.09795459659893
$block3[7,0] = 20.81648354183108
$block3[8,0] = 20.63653637581884
$block3[9,0] = 20.56054008626154
$block3[10,0] = 20.53260856928106
$block3[11,0] = 20.5385864090179
$block3[12,0] = 20.55822515495012
$block3[13,0] = 20.57036481810651
$block3[14,0] = 20.64162118768181
$block3[15,0] = 20.6868387419094
$block3[16,0] = 20.71339827634297
$block3[17,0] = 20.72248550682258
$block3[18,0] = 20.68196814244185
$block3[19,0] = 20.5524337728802
$block3[20,0] = 20.47271210185771
$block3[21,0] = 20.5148081552643
$block3[22,0] = 20.6841683858847
$block3[23,0] = 20.88792579653446
$ws.Range("O2:O25").Value = $block3
